# Attempt to collapse EVTYPE: highlight the two question paragraphs in
# yellow, and nudge the "_GoBack" bookmark (an attempted edit point) onto
# the paragraph right after them.

$d = $word.ActiveDocument

# Locate the two question paragraphs by their text so the script is not
# brittle against paragraph-index drift.
$q1 = $null
$q2 = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "Across the United States, which types of events (as indicated in the*") {
        $q1 = $p
    } elseif ($t -like "Across the United States, which types of events have the greatest economic consequences?*") {
        $q2 = $p
    }
}

# Highlight both paragraphs (text + paragraph mark) in yellow.
$q1.Range.Font.HighlightColorIndex = 7
$q2.Range.Font.HighlightColorIndex = 7

# Move the "_GoBack" bookmark to the very start of the paragraph that
# follows the two questions (the "Consider writing your report..." one).
$q2Index = $q2.Index
$target = $d.Paragraphs($q2Index + 1)
$d.Bookmarks.Add("_GoBack", $d.Range($target.Range.Start, $target.Range.Start))
